$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update existing values (Occupancy data correction) ---
$ws.Range("B3").Value = 847715

$ws.Range("C4").Value = 1329409

$ws.Range("B5").Value = 2702635
$ws.Range("C5").Value = 1314178

$ws.Range("C6").Value = 539911

$ws.Range("B7").Value = 1746137

$ws.Range("B8").Value = 2319236
$ws.Range("C8").Value = 712032

$ws.Range("C9").Value = 2829521

$ws.Range("B10").Value = 688581
$ws.Range("C10").Value = 937646

$ws.Range("B12").Value = 2041851
$ws.Range("C12").Value = 1255315

$ws.Range("B13").Value = 2961973

$ws.Range("B14").Value = 1965887
$ws.Range("C14").Value = 617421

# --- Append new rows 15-20 with fresh data, extending the table ---
$newRows = @(
    @{ Row = 15; A = 43538; B = 2007874; C = 1754423 },
    @{ Row = 16; A = 43539; B = 1712458; C = 704725 },
    @{ Row = 17; A = 43540; B = 601405;  C = 2502100 },
    @{ Row = 18; A = 43541; B = 96503;   C = 2240032 },
    @{ Row = 19; A = 43542; B = 2555699; C = 1355297 },
    @{ Row = 20; A = 43543; B = 2025443; C = 2058823 }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    # Copy the date cell above so the date number format (style) carries over
    $ws.Range("A" + ($rowNum - 1)).Copy($ws.Range("A" + $rowNum))
    $ws.Range("A" + $rowNum).Value = $r.A
    $ws.Range("B" + $rowNum).Value = $r.B
    $ws.Range("C" + $rowNum).Value = $r.C
    $ws.Range("D" + $rowNum).Formula = '=B' + $rowNum + '+C' + $rowNum + '/Hoja2!$A$2'
}

# --- Extend the table (ListObject) to cover the new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D20"))

# --- Reset the selection to A1 (matches a cleared/default selection) ---
$ws.Range("A1").Select()
